{"js": "// TIMCI Pragmatic Cluster RCT Monitoring Report \u2014 \"Minor update of the style\n// export\":\n//   1. Word's auto \"last edit position\" bookmark (`_GoBack`) moves from the\n//      end of the \"Number of children\" heading to a point inside the title\n//      text (right after \"TIMCI Pragmatic Clu\"), splitting that run in two.\n//   2. The four heading/title styles (Title, Heading1, Heading2, Heading3)\n//      get their font color pinned to a literal RGB value instead of a\n//      theme-color reference (as happens when Word/the style-export re-writes\n//      them with resolved colors).\n\nconst doc = context.document;\nconst body = doc.body;\n\n// --- 1. Relocate the \"_GoBack\" bookmark -----------------------------------\n// Remove the existing one (wherever Word last parked it) \u2026\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// \u2026 then drop a fresh one right after \"TIMCI Pragmatic Clu\" in the title.\nconst hits = body.search(\"TIMCI Pragmatic Clu\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  const splitPoint = hits.items[0].getRange(\"End\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 2. Re-point the heading/title styles at literal RGB colors -----------\nconst styles = context.document.getStyles();\nconst title = styles.getByNameOrNullObject(\"Title\");\nconst heading1 = styles.getByNameOrNullObject(\"Heading 1\");\nconst heading2 = styles.getByNameOrNullObject(\"Heading 2\");\nconst heading3 = styles.getByNameOrNullObject(\"Heading 3\");\nawait context.sync();\n\ntitle.font.color = \"#3D527A\";\nheading1.font.color = \"#3D527A\";\nheading2.font.color = \"#79BBD6\";\nheading3.font.color = \"#FC7969\";\nawait context.sync();\n", "ps1": "# TIMCI Pragmatic Cluster RCT Monitoring Report -- \"Minor update of the style\n# export\":\n#   1. Word's auto \"last edit position\" bookmark (_GoBack) moves from the end\n#      of the \"Number of children\" heading to a point inside the title text\n#      (right after \"TIMCI Pragmatic Clu\"), splitting that run in two.\n#   2. The four heading/title styles (Title, Heading1, Heading2, Heading3)\n#      get their font color pinned to a literal RGB value instead of a\n#      theme-color reference (as happens when Word/the style-export re-writes\n#      them with resolved colors).\n\n$d = $word.ActiveDocument\n\n# --- 1. Relocate the \"_GoBack\" bookmark ------------------------------------\n# Remove the existing one (wherever Word last parked it) ...\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# ... then drop a fresh one right after \"TIMCI Pragmatic Clu\" in the title.\n$splitPoint = $d.Content\n$splitPoint.Find.Execute(\"TIMCI Pragmatic Clu\") | Out-Null\nif ($splitPoint.Find.Found) {\n    $splitPoint.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $splitPoint) | Out-Null\n}\n\n# --- 2. Re-point the heading/title styles at literal RGB colors -----------\nfunction ToWordColor([string]$hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return ($b * 65536) + ($g * 256) + $r\n}\n\n$d.Styles(\"Title\").Font.Color = ToWordColor(\"3D527A\")\n$d.Styles(\"Heading 1\").Font.Color = ToWordColor(\"3D527A\")\n$d.Styles(\"Heading 2\").Font.Color = ToWordColor(\"79BBD6\")\n$d.Styles(\"Heading 3\").Font.Color = ToWordColor(\"FC7969\")\n"}
